# Refresh the COVID-19 "paises" data table (sheet "Pais") with the
# 23:26 snapshot: updated case counts for several countries, plus two
# pairs of countries that swapped rank order (and therefore row
# position) in the underlying (descending, by total cases) sort:
#   - Siria overtook Trinidad yTobago
#   - Burkina Faso overtook Nueva Zelanda
#   - Santa Lucia / Timor Oriental and Montserrat / Islas Malvinas were
#     tied and flipped order too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 23:26"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 2).Value = 6960651
$ws.Cells.Item(4, 3).Value = 34710
$ws.Cells.Item(4, 4).Value = 4214502
$ws.Cells.Item(4, 5).Value = 2542454
$ws.Cells.Item(4, 7).Value = 528
$ws.Cells.Item(4, 8).Value = 203695

# Row 5: India -> India
$ws.Cells.Item(5, 2).Value = 5398230
$ws.Cells.Item(5, 3).Value = 92755
$ws.Cells.Item(5, 4).Value = 4299724
$ws.Cells.Item(5, 5).Value = 1011732
$ws.Cells.Item(5, 7).Value = 1149
$ws.Cells.Item(5, 8).Value = 86774

# Row 6: Brasil -> Brasil
$ws.Cells.Item(6, 2).Value = 4528240
$ws.Cells.Item(6, 3).Value = 30806
$ws.Cells.Item(6, 4).Value = 3820095
$ws.Cells.Item(6, 5).Value = 571613
$ws.Cells.Item(6, 7).Value = 675
$ws.Cells.Item(6, 8).Value = 136532

# Row 29: Canada -> Canada
$ws.Cells.Item(29, 2).Value = 142774
$ws.Cells.Item(29, 3).Value = 863
$ws.Cells.Item(29, 4).Value = 124187
$ws.Cells.Item(29, 5).Value = 9376

# Row 83: Costa de Marfil -> Costa de Marfil
$ws.Cells.Item(83, 2).Value = 19269
$ws.Cells.Item(83, 3).Value = 69
$ws.Cells.Item(83, 5).Value = 757

# Row 116: Cabo Verde -> Cabo Verde
$ws.Cells.Item(116, 2).Value = 5186
$ws.Cells.Item(116, 3).Value = 45
$ws.Cells.Item(116, 4).Value = 4581
$ws.Cells.Item(116, 5).Value = 555

# Row 130: Trinidad yTobago -> Siria
$ws.Cells.Item(130, 1).Value = "Siria"
$ws.Cells.Item(130, 2).Value = 3765
$ws.Cells.Item(130, 3).Value = 34
$ws.Cells.Item(130, 4).Value = 932
$ws.Cells.Item(130, 5).Value = 2663
$ws.Cells.Item(130, 7).Value = 2
$ws.Cells.Item(130, 8).Value = 170

# Row 131: Siria -> Trinidad yTobago
$ws.Cells.Item(131, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(131, 2).Value = 3739
$ws.Cells.Item(131, 3).Value = 88
$ws.Cells.Item(131, 4).Value = 1586
$ws.Cells.Item(131, 5).Value = 2093
$ws.Cells.Item(131, 8).Value = 60

# Row 155: Nueva Zelanda -> Burkina Faso
$ws.Cells.Item(155, 1).Value = "Burkina Faso"
$ws.Cells.Item(155, 2).Value = 1816
$ws.Cells.Item(155, 3).Value = 19
$ws.Cells.Item(155, 4).Value = 1176
$ws.Cells.Item(155, 5).Value = 584
$ws.Cells.Item(155, 8).Value = 56

# Row 156: Burkina Faso -> Nueva Zelanda
$ws.Cells.Item(156, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(156, 2).Value = 1811
$ws.Cells.Item(156, 3).Value = 2
$ws.Cells.Item(156, 4).Value = 1719
$ws.Cells.Item(156, 5).Value = 67
$ws.Cells.Item(156, 8).Value = 25

# Row 157: Togo -> Togo
$ws.Cells.Item(157, 2).Value = 1659
$ws.Cells.Item(157, 3).Value = 19
$ws.Cells.Item(157, 4).Value = 1259
$ws.Cells.Item(157, 5).Value = 359

# Row 204: Timor Oriental -> Santa Lucia
$ws.Cells.Item(204, 1).Value = "Santa Lucia"

# Row 205: Santa Lucia -> Timor Oriental
$ws.Cells.Item(205, 1).Value = "Timor Oriental"

# Row 214: Islas Malvinas -> Montserrat
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

# Row 215: Montserrat -> Islas Malvinas
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0

